# Auto-generated edit script: updates crypto price/volume table
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.082.95"
$ws.Range("E2").Value = "  -4.47%  "
$ws.Range("D3").Value = "2.450.44"
$ws.Range("E3").Value = "  -7.09%  "
$ws.Range("D5").Value = "'546.75"
$ws.Range("E5").Value = "  -5.73%  "
$ws.Range("D6").Value = "'145.85"
$ws.Range("E6").Value = "  -7.12%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").Value = "'0.586"
$ws.Range("E8").Value = "  -7.31%  "
$ws.Range("D9").Value = "2.448.01"
$ws.Range("E9").Value = "  -7.11%  "
$ws.Range("D10").Value = "'0.106"
$ws.Range("E10").Value = "  -10.96%  "
$ws.Range("B11").Value = "TRON"
$ws.Range("C11").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D11").Value = "'0.154"
$ws.Range("E11").Value = "  -2.04%  "
$ws.Range("B12").Value = "Toncoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D12").Value = "'5.40"
$ws.Range("E12").Value = "  -7.41%  "
$ws.Range("D13").Value = "'0.350"
$ws.Range("E13").Value = "  -9.23%  "
$ws.Range("D14").Value = "'25.86"
$ws.Range("E14").Value = "  -10.36%  "
$ws.Range("D15").Value = "2.893.07"
$ws.Range("E15").Value = "  -7.06%  "
$ws.Range("D16").Value = "'0.0000165"
$ws.Range("E16").Value = "  -10.94%  "
$ws.Range("D17").Value = "61.032.05"
$ws.Range("E17").Value = "  -4.43%  "
$ws.Range("D18").Value = "2.450.54"
$ws.Range("E18").Value = "  -7.29%  "
$ws.Range("D19").Value = "'11.04"
$ws.Range("E19").Value = "  -9.62%  "
$ws.Range("D20").Value = "'7.02"
$ws.Range("E20").Value = "  -9.70%  "
$ws.Range("D21").Value = "'4.14"
$ws.Range("E21").Value = "  -8.76%  "
$ws.Range("D22").Value = "'316.15"
$ws.Range("E22").Value = "  -8.26%  "
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").Value = "'1.85"
$ws.Range("E24").Value = "  -2.48%  "
$ws.Range("D25").Value = "'63.64"
$ws.Range("E25").Value = "  -6.99%  "
$ws.Range("B26").Value = "WrappedeETH"
$ws.Range("C26").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D26").Value = "2.581.46"
$ws.Range("E26").Value = "  -6.79%  "
$ws.Range("B27").Value = "Bittensor"
$ws.Range("C27").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D27").Value = "'545.21"
$ws.Range("E27").Value = "  -6.85%  "
$ws.Range("B28").Value = "PEPE"
$ws.Range("C28").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D28").Value = "0.0₃0956"
$ws.Range("E28").Value = "  -15.58%  "
$ws.Range("E29").Value = "  +0.42%  "
$ws.Range("D30").Value = "'1.45"
$ws.Range("E30").Value = "  -12.39%  "
$ws.Range("D31").Value = "'8.17"
$ws.Range("E31").Value = "  -11.97%  "
$ws.Range("D32").Value = "'7.52"
$ws.Range("E32").Value = "  -9.99%  "
$ws.Range("D33").Value = "'0.145"
$ws.Range("E33").Value = "  -9.84%  "
$ws.Range("D34").Value = "'1.88"
$ws.Range("E34").Value = "  -8.61%  "
$ws.Range("D35").Value = "'1.57"
$ws.Range("E35").Value = "  -10.22%  "
$ws.Range("D36").Value = "'5.80"
$ws.Range("E36").Value = "  -13.32%  "
$ws.Range("D37").Value = "'0.998"
$ws.Range("E37").Value = "  -0.11%  "
$ws.Range("D38").Value = "'4.78"
$ws.Range("E38").Value = "  -13.18%  "
$ws.Range("D39").Value = "'0.377"
$ws.Range("E39").Value = "  -6.93%  "
$ws.Range("D40").Value = "'18.31"
$ws.Range("E40").Value = "  -7.61%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "'1.75"
$ws.Range("E41").Value = "  -9.09%  "
$ws.Range("B42").Value = "Monero"
$ws.Range("C42").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D42").Value = "'140.75"
$ws.Range("E42").Value = "  -8.88%  "
$ws.Range("E43").Value = "  +0.08%  "
$ws.Range("D44").Value = "'40.29"
$ws.Range("E44").Value = "  -4.46%  "
$ws.Range("D45").Value = "'2.31"
$ws.Range("E45").Value = "  -10.97%  "
$ws.Range("D46").Value = "'145.71"
$ws.Range("E46").Value = "  -10.53%  "
$ws.Range("D47").Value = "'3.56"
$ws.Range("E47").Value = "  -9.47%  "
$ws.Range("D48").Value = "'21.26"
$ws.Range("E48").Value = "  -12.50%  "
$ws.Range("D49").Value = "'0.0531"
$ws.Range("E49").Value = "  -9.94%  "
$ws.Range("D50").Value = "'0.585"
$ws.Range("E50").Value = "  -7.77%  "
$ws.Range("D51").Value = "'0.0929"
$ws.Range("E51").Value = "  -7.75%  "
